$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Group 1: rows 552-576 -> Maria Alejandra Giraldo Franco ("MaAl")
for ($r = 552; $r -le 576; $r++) {
    $ws.Range("A$r").Value = 30000120847
    $ws.Range("H$r").Value = "30000120847_MaAl"
}

# Group 2: rows 727-751 -> Maria Paulina Jaramillo Martinez ("MaPa")
for ($r = 727; $r -le 751; $r++) {
    $ws.Range("A$r").Value = 30000096203
    $ws.Range("H$r").Value = "30000096203_MaPa"
}
